$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Split column A's width definition off from column B's (they used to share
# one <col> entry) by touching column A's Hidden flag - a no-op visually,
# but it forces the engine to record column A's 30.7109375 width on its own.
$ws.Columns("A").Hidden = $false

# Insert two new rows before row 13 (Docentes responsaveis content rows),
# pushing the remaining rows (old 13-23) down to (15-25).
$ws.Rows("13:14").Insert()

# The insert copied formatting (incl. label cell in column A) from the row
# above into the two new rows; the target layout has no A-column cell on
# these two rows, only B/C content cells, so drop the leftover A cells.
$ws.Range("A13:A14").Clear()

# Give the new B/C cells on rows 13-14 the same formatting used elsewhere
# for the B/C content columns (wrap text, vertical-top, red font in C).
$ws.Range("B10:C10").Copy()
$ws.Range("B13:C14").PasteSpecial(-4122)

# Row 13: Docentes responsaveis - first professor
$ws.Range("B13").Value = "3444370 - Rita de Cássia Lacerda Brambilla Rodrigues"
$ws.Range("C13").Value = "3444370 - Rita de Cássia Lacerda Brambilla Rodrigues"

# Row 14: Docentes responsaveis - second professor (new)
$ws.Range("B14").Value = "5817181 - Valdeir Arantes"
$ws.Range("C14").Value = "5817181 - Valdeir Arantes"

# Row 10: Objetivos - replace with the real objectives text
$ws.Range("B10").Value = "Capacitar os alunos a utilizarem conceitos e ferramentas de matemática financeira para avaliações econômicas de projetos de engenharia"
$ws.Range("C10").Value = "Capacitar os alunos a utilizarem conceitos e ferramentas de matemática financeira para avaliações econômicas de projetos de engenharia"

# Row 15: Programa resumido - replace with the real short-syllabus text
$ws.Range("B15").Value = "Introdução à Engenharia Econômica; estimativa de custos; juros; fluxo de caixa; depreciação; comparação entre alternativa de investimentos; ponto de equilíbrio; planilha eletrônica aplicado à Engenharia Econômica."
$ws.Range("C15").Value = "Introdução à Engenharia Econômica; estimativa de custos; juros; fluxo de caixa; depreciação; comparação entre alternativa de investimentos; ponto de equilíbrio; planilha eletrônica aplicado à Engenharia Econômica."

# Row 17: Programa - replace with the full syllabus text
$ws.Range("B17").Value = "1.INTRODUÇÃO À ENGENHARIA ECONÔMICA: necessidades de uma análise econômica de projeto (engenharia econômica como uma ferramenta de análise de tomada de decisão); 2.ESTIMATIVA DE CUTOS: estimativa de investimento de capital (tipos de estimativas de custo de capital; índice de custos mais comuns; métodos para estimar investimento de capital); estimativa de custo de equipamentos; estimativa de custo de produção;3.JUROS: variável tempo (juros simples; juros compostos); taxa efetiva, nominal e equivalente; relações de equivalência. 4.FLUXO DE CAIXA: diagrama de fluxo de caixa; equivalência de fluxo de caixa; séries uniforme e em gradientes; elaboração de fluxo de caixa5.Depreciação: métodos de depreciação;6.COMPARAÇÃO ENTRE ALTERNATIVAS DE INVESTIMENTO: critério de Rentabilidade – método do Valor Anual Uniforme Equivalente (VAU); método do Valor Presente (VP); método da Taxa interna de Retorno (TIR); método da Taxa Interna de Retorno Modificada (TIRM); renovação e substituição de equipamentos; período de retorno;7.Ponto de Equilíbrio de um projeto (break-even-point);8.Planilhas eletrônicas aplicado à análises de Engenharia Econômica."
$ws.Range("C17").Value = "1.INTRODUÇÃO À ENGENHARIA ECONÔMICA: necessidades de uma análise econômica de projeto (engenharia econômica como uma ferramenta de análise de tomada de decisão); 2.ESTIMATIVA DE CUTOS: estimativa de investimento de capital (tipos de estimativas de custo de capital; índice de custos mais comuns; métodos para estimar investimento de capital); estimativa de custo de equipamentos; estimativa de custo de produção;3.JUROS: variável tempo (juros simples; juros compostos); taxa efetiva, nominal e equivalente; relações de equivalência. 4.FLUXO DE CAIXA: diagrama de fluxo de caixa; equivalência de fluxo de caixa; séries uniforme e em gradientes; elaboração de fluxo de caixa5.Depreciação: métodos de depreciação;6.COMPARAÇÃO ENTRE ALTERNATIVAS DE INVESTIMENTO: critério de Rentabilidade – método do Valor Anual Uniforme Equivalente (VAU); método do Valor Presente (VP); método da Taxa interna de Retorno (TIR); método da Taxa Interna de Retorno Modificada (TIRM); renovação e substituição de equipamentos; período de retorno;7.Ponto de Equilíbrio de um projeto (break-even-point);8.Planilhas eletrônicas aplicado à análises de Engenharia Econômica."

# Row 20: Metodo - replace with the real method text
$ws.Range("B20").Value = "O método utilizado tem por fundamento a Aprendizagem Baseada em Projetos (PBL) que visa desenvolver as competências técnicas relativas ao tema do projeto, bem como competências transversais, tais como: aprender a aprender, trabalho em equipe, relacionamento interpessoal, aspectos de liderança e capacidade de comunicação, dentre outras; exercícios individuais realizados no decorrer da disciplina. Para os projetos, os alunos serão divididos em grupos que desenvolverão um projeto durante o semestre relacionado a um tema de Engenharia de Econômica na área de Engenharia Bioquímica, similar ao que eles irão encontrar na vida real no efetivo exercício de sua profissão."
$ws.Range("C20").Value = "O método utilizado tem por fundamento a Aprendizagem Baseada em Projetos (PBL) que visa desenvolver as competências técnicas relativas ao tema do projeto, bem como competências transversais, tais como: aprender a aprender, trabalho em equipe, relacionamento interpessoal, aspectos de liderança e capacidade de comunicação, dentre outras; exercícios individuais realizados no decorrer da disciplina. Para os projetos, os alunos serão divididos em grupos que desenvolverão um projeto durante o semestre relacionado a um tema de Engenharia de Econômica na área de Engenharia Bioquímica, similar ao que eles irão encontrar na vida real no efetivo exercício de sua profissão."

# Row 21: Criterio - replace with the real criterion text
$ws.Range("B21").Value = "A nota do projeto (NP) será individual e será a média ponderada de componentes do projeto, tais como: Projeto Preliminar, Projeto Final, envolvimento do aluno com o projeto, avaliação dos Pares, apresentação de trabalhos, dentre outros;A nota dos exercícios (NE) será a média aritmética dos exercícios individuais;NF = (0,8NP + 0,2NE)"
$ws.Range("C21").Value = "A nota do projeto (NP) será individual e será a média ponderada de componentes do projeto, tais como: Projeto Preliminar, Projeto Final, envolvimento do aluno com o projeto, avaliação dos Pares, apresentação de trabalhos, dentre outros;A nota dos exercícios (NE) será a média aritmética dos exercícios individuais;NF = (0,8NP + 0,2NE)"

# Row 22: Norma de recuperacao - replace with the real recovery-rule text
$ws.Range("B22").Value = "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2n"
$ws.Range("C22").Value = "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2n"

# Row 23: Bibliografia - replace with the real bibliography text
$ws.Range("B23").Value = "1.Blank L. e Tarquin A. Engenharia Econômica. 6ed. Mc GrawHill, 2007.2.Fraser N. M., Jewkes, Elizabeth M. Engineering Economics: Financial decision making for engineers. 5th Edition. Editora Pearson, 2013. 3.Samanez, Carlos Patrício. Engenharia Econômica. São Paulo. Editora Prentice Hall, 2009. 4.Hirschfeld, Henrique. Engenharia Econômica e análise de custos. 7ed. Editora Atalas, 2007. 5.Ehrlich, Pierre Jacques. Engenharia Econômica. São Paulo. Editora Atlas, 2005 6.Peter, MS. And Timmerhaus, K.D. Plant Design and economic for chemical engineers"
$ws.Range("C23").Value = "1.Blank L. e Tarquin A. Engenharia Econômica. 6ed. Mc GrawHill, 2007.2.Fraser N. M., Jewkes, Elizabeth M. Engineering Economics: Financial decision making for engineers. 5th Edition. Editora Pearson, 2013. 3.Samanez, Carlos Patrício. Engenharia Econômica. São Paulo. Editora Prentice Hall, 2009. 4.Hirschfeld, Henrique. Engenharia Econômica e análise de custos. 7ed. Editora Atalas, 2007. 5.Ehrlich, Pierre Jacques. Engenharia Econômica. São Paulo. Editora Atlas, 2005 6.Peter, MS. And Timmerhaus, K.D. Plant Design and economic for chemical engineers"

Write-Host "done"
